$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"

$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been executed successfully. Both `"Oppenheimer`" and `"Barbie`" will be showcased as planned.`n"

$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to both `"Barbie`" and `"Oppenheimer.`"`n"

$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on the movie.`n"
$ws.Range("D5").Value = "no_decision, "

$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.`n"

$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded successfully.`n"

$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"

$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding which movie will be shown on Friday.`n"
$ws.Range("D9").Value = "no_decision, "

$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D10").Value = "no_decision, "

$ws.Range("C11").Value = "MSG: None`n`nMSG: The function has been called, indicating that no decision about Friday’s movie has been reached.`n"

$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made, as the committee did not reach an agreement.`n"
$ws.Range("D12").Value = "no_decision, "

$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Oppenheimer`" will be shown on Friday.`n"

$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday has not been reached, and no definitive choice has been made.`n"
$ws.Range("D14").Value = "no_decision, "

$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"

$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Oppenheimer`" will be acquired for the screening on Friday.`n"

$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision process has concluded without a choice of movie for Friday, resulting in no decision being made.`n"
$ws.Range("D17").Value = "no_decision, "

$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision-making process concluded without a consensus on which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("D18").Value = "no_decision, "

$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"

$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision-making process did not result in a clear choice for a movie to show on Friday, so the result is that no decision was made.`n"
$ws.Range("D20").Value = "no_decision, "

$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected.`n"
$ws.Range("D21").Value = "no_decision, "

$ws.Range("C22").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D22").Value = "both_movies, "

$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision-making process concluded without a clear movie selection for Friday, resulting in no decision being made.`n"
$ws.Range("D23").Value = "no_decision, "

$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision process indicates no final agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("D24").Value = "no_decision, "

$ws.Range("C25").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D25").Value = "both_movies, "

$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("D26").Value = "no_decision, "

$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded, and the outcome is that there is no decision regarding the movie for Friday.`n"
$ws.Range("D27").Value = "no_decision, "

$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.`n"

$ws.Range("C29").Value = "MSG: None`n`nMSG: I have recorded the decision as there was no consensus on a movie for Friday.`n"
$ws.Range("D29").Value = "no_decision, "

$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made regarding the movie to be shown on Friday.`n"
$ws.Range("D30").Value = "no_decision, "

$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded successfully.`n"

$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D32").Value = "no_decision, "

$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has not been made.`n"

$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"

$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded as no definitive choice for Friday's movie has been made.`n"
$ws.Range("D35").Value = "no_decision, "

$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D36").Value = "no_decision, "

$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision about the movie for Friday was not reached, so no rights will be acquired.`n"
$ws.Range("D37").Value = "no_decision, "

$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as there was no consensus on the movie to be shown on Friday.`n"
$ws.Range("D38").Value = "no_decision, "

$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision about Friday's movie can be made.`n"
$ws.Range("D39").Value = "no_decision, "

$ws.Range("C40").Value = "MSG: None`n`nMSG: The committee has not reached a decision regarding which movie to show on Friday. As a result, the no_decision function was called, indicating that no agreement was made.`n"
$ws.Range("D40").Value = "no_decision, "

$ws.Range("C41").Value = "MSG: None`n`nMSG: I have recorded the decision that no movie has been selected.`n"
$ws.Range("D41").Value = "no_decision, "

$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision process has concluded without a definitive choice for Friday's movie.`n"
$ws.Range("D42").Value = "no_decision, "

$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"

$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Oppenheimer`" have been acquired.`n"

$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`".`n"
$ws.Range("D45").Value = "no_decision, "

$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision-making process has concluded without a definitive choice for Friday's movie, resulting in no film being selected for acquisition.`n"
$ws.Range("D46").Value = "no_decision, "

$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been made to show `"Barbie`" on Friday.`n"

$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has been recorded as `"no decision.`"`n"
$ws.Range("D48").Value = "no_decision, "

$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("D49").Value = "no_decision, "

$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"

$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"

$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision about the movie to be shown on Friday has ended without a selection.`n"
$ws.Range("D52").Value = "no_decision, "

$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been made to acquire rights for the movie `"Barbie`" for Friday's showing.`n"

$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision-making process did not result in an agreement on a movie to show on Friday, leading to no decision being reached.`n"
$ws.Range("D54").Value = "no_decision, "

$ws.Range("C55").Value = "MSG: None`n`nMSG: No decision was made regarding the movie for Friday.`n"
$ws.Range("D55").Value = "no_decision, "

$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded, and there will be no movie shown on Friday.`n"
$ws.Range("D56").Value = "no_decision, "

$ws.Range("C57").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D57").Value = "both_movies, "

$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D58").Value = "no_decision, "

$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been recorded as no clear agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("D59").Value = "no_decision, "

$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D60").Value = "no_decision, "

$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has concluded without a clear choice.`n"
$ws.Range("D61").Value = "no_decision, "

$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to the movie `"Barbie`".`n"

$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday remains unresolved, resulting in no agreement being reached.`n"
$ws.Range("D63").Value = "no_decision, "

$ws.Range("C64").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired. The committee has decided to show `"Oppenheimer`" in its entirety, along with the inclusion of `"Barbie.`"`n"
$ws.Range("D64").Value = "both_movies, "

$ws.Range("C65").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no decision.`"`n"
$ws.Range("D65").Value = "no_decision, "

$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" was selected for the assembly.`n"
